# Apply updated transition-probability values to the team-specific matrix
# sheet. These numbers reflect additional simulated games, so many cells
# that were previously 0 (no observed transition) now hold the observed
# frequency for that starting/ending state pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = 1

    "P3" = 1

    "P4" = 1

    "J6" = 0.2857142857142857
    "O6" = 0.1428571428571428
    "R6" = 0.2857142857142857
    "S6" = 0.2857142857142857

    "D7" = 0.3333333333333333
    "S7" = 0.6666666666666666

    "F8" = 0.1
    "Q8" = 0.3
    "R8" = 0.1
    "S8" = 0.5

    "J9" = 0.5
    "Q9" = 0.25
    "S9" = 0.25

    "B10" = 0.03389830508474576
    "D10" = 0.03389830508474576
    "F10" = 0.05084745762711865
    "J10" = 0.03389830508474576
    "O10" = 0.01694915254237288
    "Q10" = 0.2542372881355932
    "R10" = 0.1694915254237288
    "S10" = 0.4067796610169492

    "G11" = 0.1428571428571428
    "J11" = 0.2857142857142857
    "K11" = 0.2857142857142857
    "L11" = 0.2857142857142857

    "G12" = 1

    "H15" = 0.1111111111111111
    "J15" = 0.5555555555555556
    "O15" = 0.1111111111111111
    "S15" = 0.2222222222222222

    "H16" = 0.6
    "J16" = 0.2
    "O16" = 0.2

    "H17" = 0.1578947368421053
    "I17" = 0.05263157894736842
    "J17" = 0.6842105263157895
    "K17" = 0.05263157894736842
    "S17" = 0.05263157894736842

    "H18" = 0.07692307692307693
    "I18" = 0.07692307692307693
    "J18" = 0.7692307692307693
    "O18" = 0.07692307692307693

    "F19" = 0.05
    "H19" = 0.075
    "I19" = 0.05
    "J19" = 0.55
    "K19" = 0.1
    "O19" = 0.1
    "S19" = 0.1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
